# Generate Report for Handoff
# Updates the "Latest Handoff Date"/"Latest Handoff Datetime" timestamps
# for the rows that just had a new handoff generated (all rows except the
# "Handed back: in sync with en-US" rows and the "In Translation" row).

$wb = $excel.ActiveWorkbook

$rows = @(4, 6, 7, 8, 9, 10)

# "Overview" sheet — column D is "Latest Handoff Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 4).Value = "2016-03-22 07:17:54"
}

# "zh-cn" sheet — column E is "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "2016-03-22 07:17:46"
}

# "de-de" sheet — column E is "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "2016-03-22 07:17:54"
}
